$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell B5: 96603 -> 96604
$ws.Range("B5").Value = 96604

# Add new row 11 with observation data
$ws.Range("A11").Value = 131223036
$ws.Range("B11").Value = 57881

$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 100049
$ws.Range("F11").Value = "Spillkråka"
$ws.Range("G11").Value = "Dryocopus martius"
$ws.Range("H11").Value = "(Linnaeus, 1758)"

# I11 holds the text "1" (not a numeric value) in the source data
$ws.Range("I11").NumberFormat = "@"
$ws.Range("I11").Value = "1"
$ws.Range("I11").Style = "Normal"

# K11 / L11 are present but blank in the source data
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""

$ws.Range("M11").Value = "spel/sång"
$ws.Range("N11").Value = "passiv ljudinspelning"

$ws.Range("P11").Value = "Korseberget, Boh"
$ws.Range("Q11").Value = 311091
$ws.Range("R11").Value = 6410588
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = "Västra Götaland"
$ws.Range("U11").Value = "Kungälv"
$ws.Range("V11").Value = "Bohuslän"
$ws.Range("W11").Value = "Harestad"

# Dates/times must stay plain text, not auto-converted to date serials
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y11").Value = "2026-02-11"
$ws.Range("Y11").Style = "Normal"

$ws.Range("Z11").NumberFormat = "@"
$ws.Range("Z11").Value = "13:30"
$ws.Range("Z11").Style = "Normal"

$ws.Range("AA11").NumberFormat = "@"
$ws.Range("AA11").Value = "2026-02-11"
$ws.Range("AA11").Style = "Normal"

$ws.Range("AB11").NumberFormat = "@"
$ws.Range("AB11").Value = "13:30"
$ws.Range("AB11").Style = "Normal"

$ws.Range("AC11").Value = "Inspelad i fält på platsen med en Audiomoth  inspelningsapparat"

$ws.Range("AD11").Value = $False
$ws.Range("AE11").Value = $False
$ws.Range("AG11").Value = $False

$ws.Range("AT11").Value = ""

$ws.Range("AW11").Value = "Linus Lundin"
$ws.Range("AX11").Value = "Linus Lundin"

$ws.Range("AY11").Value = ""
